# Add a new search criterion property definition row for "user.parent.id"
# (parent of a user), inserted right before the existing row that defines
# "user.decrypt" (row 269), pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new row at position 269 (existing rows 269.. shift to 270..)
[void]$ws.Rows.Item(269).Insert()
$ws.Rows.Item(269).RowHeight = 16.5

# Populate the new row's cells
#   A = entity            -> USER_DB
#   B = property           -> user.parent.id
#   C = value_type          -> LONG
#   K = picker             -> USER_DB
#   L = name_l10n_key        -> user.parent.id
#   M = valid_restrictions     -> EQ, NE, IS_NULL
$ws.Cells.Item(269, 1).Value = "USER_DB"
$ws.Cells.Item(269, 2).Value = "user.parent.id"
$ws.Cells.Item(269, 3).Value = "LONG"
$ws.Cells.Item(269, 11).Value = "USER_DB"
$ws.Cells.Item(269, 12).Value = "user.parent.id"
$ws.Cells.Item(269, 13).Value = "EQ, NE, IS_NULL"

# Update the view state to reflect where the edit was made
[void]$ws.Range("A269").Select()
$excel.ActiveWindow.ScrollRow = 261
$excel.ActiveWindow.ScrollColumn = 1
